$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Map Construction" block in columns N:P.
# Cells are written in the same order the strings were first introduced
# so the shared-strings table comes out in the expected order.

$ws.Range("N1").Value = "Map Construction"

$ws.Range("N3").Value = "Format"

$ws.Range("N5").Value = "Width"
$ws.Range("N6").Value = "Height"
$ws.Range("N7").Value = "Start X Position"
$ws.Range("N8").Value = "Start Y Position"
$ws.Range("N9").Value = "NPC Count"
$ws.Range("N10").Value = "Map data"

$ws.Range("O4").Value = "n/a"
$ws.Range("P4").Value = "n/a"

$ws.Range("O3").Value = "Min"
$ws.Range("P3").Value = "Max"

$ws.Range("N4").Value = "Name (alpha-numeric)"

$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 50

$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 50

$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 49

$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 49

$ws.Range("O9").Value = 0
$ws.Range("P9").Formula = "=(50*50)-1"

$ws.Range("O10").Value = "n/a"
$ws.Range("P10").Value = "n/a"

# Columns O:P mirror the left-aligned "style 2" formatting already used by
# the other data columns (D, I) in this sheet.
$ws.Range("O3:P10").HorizontalAlignment = -4131

# Column N is sized to fit its longest label; O:P stay at the default width.
$ws.Columns.Item(14).ColumnWidth = 21.42578125
$ws.Columns.Item(15).ColumnWidth = 9.140625
$ws.Columns.Item(16).ColumnWidth = 9.140625

# Scroll the view over to show the new columns and match the author's
# final selection.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("M8").Select()
